# Add a new "consequents_length" column (H) to the rule-extraction sheet,
# mirroring the existing "antecedents_length" column (G):
#   - H1 header text, same bold/border/centered style as the other headers
#   - H2 data value (1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in H1.
$ws.Range("H1").Value = "consequents_length"

# Give H1 the same formatting as the neighbouring header cell (G1) - bold
# font, thin border, centered alignment - by copying just the formatting
# (xlPasteFormats = -4122) rather than the value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell in H2.
$ws.Range("H2").Value = 1
